$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update migration wave date in title cell
$ws.Range("A1").Value = "Hotcarding Spreadsheet - Migration Wave 10/22/2002"

# Add new data row 5, mirroring the pattern of row 4
# A5 holds a date-looking string; force text format so it is not
# auto-converted into a date serial number.
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2025-10-01"
$ws.Range("B5").Value = "FinanceCorp"
$ws.Range("C5").Value = "ENT12375"
$ws.Range("D5").Value = "FISB"
$ws.Range("E5").Value = "LegacyPay"
$ws.Range("F5").Value = "PaymentsOne Debit"
$ws.Range("G5").Value = "Basic"
$ws.Range("H5").Value = "Offshore"
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = "Yes"

# Clear the explicit per-cell style the engine stamped on write so the
# new cells rely on the column default style, matching row 4's pattern.
$ws.Range("A5:J5").Style = "Normal"
